# Add data for 2022-08-04
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-07-27"

# Update the header label in I1 to reflect the new "through" date
$ws.Range("I1").Value = "2022 (through 07-27)"

# Update the July total (row 8) and the Total row (row 14)
$ws.Range("I8").Value = 151
$ws.Range("I14").Value = 957
